$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

Set-CellText $ws "D2" "69.335.66"
Set-CellText $ws "E2" "  +0.06%  "
Set-CellText $ws "D3" "3.658.47"
Set-CellText $ws "E3" "  -0.75%  "
Set-CellText $ws "E4" "  +0.03%  "
Set-CellText $ws "D5" "641.29"
Set-CellText $ws "E5" "  -5.60%  "
Set-CellText $ws "D6" "158.19"
Set-CellText $ws "E6" "  -0.53%  "
Set-CellText $ws "E7" "  -0.03%  "
Set-CellText $ws "D8" "0.495"
Set-CellText $ws "E8" "  +0.31%  "
Set-CellText $ws "E9" "  -1.26%  "
Set-CellText $ws "D10" "7.00"
Set-CellText $ws "E10" "  -1.67%  "
Set-CellText $ws "D11" "0.436"
Set-CellText $ws "E11" "  -0.39%  "
Set-CellText $ws "E12" "  -1.50%  "
Set-CellText $ws "D13" "4.281.73"
Set-CellText $ws "E13" "  -0.61%  "
Set-CellText $ws "D14" "32.12"
Set-CellText $ws "E14" "  -1.00%  "
Set-CellText $ws "D15" "3.690.48"
Set-CellText $ws "E15" "  +0.14%  "
Set-CellText $ws "D16" "69.337.26"
Set-CellText $ws "E16" "  +0.07%  "
Set-CellText $ws "E17" "  +1.11%  "
Set-CellText $ws "D18" "15.87"
Set-CellText $ws "E18" "  -0.88%  "
Set-CellText $ws "D19" "6.38"
Set-CellText $ws "E19" "  -0.62%  "
Set-CellText $ws "D20" "463.69"
Set-CellText $ws "E20" "  -0.82%  "
Set-CellText $ws "D21" "9.66"
Set-CellText $ws "E21" "  -3.11%  "
Set-CellText $ws "D22" "0.638"
Set-CellText $ws "E22" "  -2.12%  "
Set-CellText $ws "D23" "79.38"
Set-CellText $ws "E23" "  -0.51%  "
Set-CellText $ws "D24" "3.811.14"
Set-CellText $ws "E24" "  -0.61%  "
Set-CellText $ws "E25" "  -0.01%  "
Set-CellText $ws "D26" "0.0000123"
Set-CellText $ws "E26" "  +0.06%  "
Set-CellText $ws "D27" "10.67"
Set-CellText $ws "E27" "  -2.39%  "
Set-CellText $ws "D28" "8.81"
Set-CellText $ws "E28" "  -3.57%  "
Set-CellText $ws "D29" "2.58"
Set-CellText $ws "E29" "  -3.49%  "
Set-CellText $ws "D31" "0.999"
Set-CellText $ws "E31" "  -0.05%  "
Set-CellText $ws "E32" "  -0.85%  "
Set-CellText $ws "D33" "26.41"
Set-CellText $ws "E33" "  -1.88%  "
Set-CellText $ws "D34" "6.38"
Set-CellText $ws "E34" "  -3.69%  "
Set-CellText $ws "D35" "3.652.18"
Set-CellText $ws "E35" "  -0.61%  "
Set-CellText $ws "E36" "  +1.63%  "
Set-CellText $ws "D37" "8.25"
Set-CellText $ws "E37" "  -0.29%  "
Set-CellText $ws "E38" "  +0.00%  "
Set-CellText $ws "D39" "5.86"
Set-CellText $ws "E39" "  -5.49%  "
Set-CellText $ws "D40" "178.22"
Set-CellText $ws "E40" "  +4.08%  "
Set-CellText $ws "D41" "1.00"
Set-CellText $ws "E41" "  +0.03%  "
Set-CellText $ws "D42" "0.0889"
Set-CellText $ws "E42" "  -1.97%  "
Set-CellText $ws "D43" "2.16"
Set-CellText $ws "E43" "  -4.37%  "
Set-CellText $ws "D44" "0.924"
Set-CellText $ws "E44" "  -2.00%  "
Set-CellText $ws "D45" "46.50"
Set-CellText $ws "E45" "  -2.32%  "
Set-CellText $ws "D46" "2.66"
Set-CellText $ws "E46" "  -0.89%  "
Set-CellText $ws "B47" "ONDO"
Set-CellText $ws "C47" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-CellText $ws "D47" "1.25"
Set-CellText $ws "E47" "  -3.17%  "
Set-CellText $ws "B48" "InjectiveProtocol"
Set-CellText $ws "C48" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-CellText $ws "D48" "26.63"
Set-CellText $ws "E48" "  -6.15%  "
Set-CellText $ws "D49" "7.74"
Set-CellText $ws "E49" "  -0.52%  "
Set-CellText $ws "B50" "SuiNetwork"
Set-CellText $ws "C50" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-CellText $ws "D50" "1.04"
Set-CellText $ws "E50" "  -6.64%  "
Set-CellText $ws "B51" "FLOKI"
Set-CellText $ws "C51" "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-CellText $ws "D51" "0.000259"
Set-CellText $ws "E51" "  -6.05%  "
